$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A5").NumberFormat = "@"

$ws.Range("A2").Value = "5171429829879"
$ws.Range("A3").Value = "4715172468769"
$ws.Range("A4").Value = "2483774646197"
$ws.Range("A5").Value = "5297607499672"
